$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: row index 0..29 for rows 2..31 ---
for ($i = 0; $i -le 29; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $i
}

# New rows 30 & 31 need the same "index column" look (bold, thin border,
# centered/top) that the pre-existing A2:A29 cells already carry.
$ws.Range("A30:A31").Font.Bold = $true
$ws.Range("A30:A31").HorizontalAlignment = -4108
$ws.Range("A30:A31").VerticalAlignment = -4160
$ws.Range("A30:A31").Borders.LineStyle = 1

# --- Column B: HKL / ring labels for rows 2..31 ---
$labels = @("HKL","Spiral5","Holden","Rizzie Spiral","RotRing OmegaMax-90","Equal Angle","Tilt Rotate","CLR","Rizzie Hex","Matthies Hex","Tilt Rotate_Partial","RotRing OmegaMax-60","Equal Angle_Partial","Rizzie Hex_Partial","ND Single","RD Single","TD Single","Morris Single","Ring Perpendicular to ND","Ring Perpendicular to RD","Ring Perpendicular to TD","OffsetFTD","OffsetATD","OffsetF45","OffsetA45","OffsetFRD","OffsetARD","Gaussian Quadrature","Michael-CCHex","Michael-SNHex")
for ($i = 0; $i -le 29; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $labels[$i]
}

# --- Row 2 header labels, columns C..T ---
$headerLabels = @("[1, 1, 0]","[3, 1, 0]","[2, 2, 2]","[3, 2, 1]","[2, 1, 1]","[4, 0, 0]","[2, 0, 0]","[2, 2, 0]","1Pair-A","1Pair-B","2Pairs-A","2Pairs-B","3Pairs-A","3Pairs-B","3Pairs-C","4Pairs","5A4F","MaxUnique")
for ($j = 0; $j -le 17; $j++) {
    $ws.Cells.Item(2, $j + 3).Value = $headerLabels[$j]
}

# --- Data block rows 3..31, columns C..T ---
$dataRows = @(
  ,(1.120079853657215,0.9961329895178055,1.077861162212376,0.990415576587013,1.019986702077031,1.079326715570113,1.079326715570113,1.120079853657215,1.120079853657215,1.019986702077031,1.049656708823572,1.049656708823572,1.031815469054983,1.073131090434786,1.073131090434786,1.084868281240393,1.084868281240393,1.047300499936925)
  ,(0.8230609169199407,0.9469505479934185,1.927691841862396,0.9337896172281688,1.289668957223146,2.225531357924238,2.225531357924238,0.8230609169199407,0.8230609169199407,1.289668957223146,1.757600157573692,1.757600157573692,1.487383621046934,1.446087077355775,1.446087077355775,1.290330537246816,1.290330537246816,1.357782206525218)
  ,(5.818668321352504,0.3019432256074346,2.154500460863252,1.724601467711075,1.530852014550058,0.2319101419496165,0.2319101419496165,5.818668321352504,5.818668321352504,1.530852014550058,0.8813810782498371,0.8813810782498371,0.6882351273690364,2.527143492617393,2.527143492617393,3.350024699801171,3.350024699801171,1.96041260533899)
  ,(1.199068779750939,0.8032113772738989,1.714103526287078,0.9801024240181715,1.114660730763665,0.904415196490881,0.904415196490881,1.199068779750939,1.199068779750939,1.114660730763665,1.009537963627273,1.009537963627273,0.9407624348428151,1.072714902335162,1.072714902335162,1.104303371689106,1.104303371689106,1.119260339097439)
  ,(1.097655098868639,0.9266058341850941,2.612722075937493,0.9299062723889632,1.512226357570853,3.155020307254655,3.155020307254655,1.097655098868639,1.097655098868639,1.512226357570853,2.333623332412754,2.333623332412754,1.8646174996702,1.921633921231382,1.921633921231382,1.715639215640696,1.715639215640696,1.705689324367616)
  ,(0.4428852883562573,0.4686054254572429,16.06109719942996,0.4692314428937391,5.682126891637458,21.24948437995802,21.24948437995802,0.4428852883562573,0.4428852883562573,5.682126891637458,13.46580563579774,13.46580563579774,9.133405565684241,9.124832186650579,9.124832186650579,6.954345462076999,6.954345462076999,7.39557177128878)
  ,(1.240254493033031,1.028613338946695,1.041821829527088,1.031618886503664,1.039614554122016,1.099851321006688,1.099851321006688,1.240254493033031,1.240254493033031,1.039614554122016,1.069732937564352,1.069732937564352,1.0560264046918,1.126573456053912,1.126573456053912,1.154993715298692,1.154993715298692,1.080295737189864)
  ,(1.137781263764229,0.9884318269651895,1.02698613553563,1.004898742426084,1.011296751049787,1.08997083278909,1.08997083278909,1.137781263764229,1.137781263764229,1.011296751049787,1.050633791919439,1.050633791919439,1.029899803601356,1.079682949201035,1.079682949201035,1.094207527841834,1.094207527841834,1.043227592088335)
  ,(1.604907144245246,1.078123803187341,1.208034769188831,1.092598002071044,1.117738869766187,1.405099968017963,1.405099968017963,1.604907144245246,1.604907144245246,1.117738869766187,1.261419418892075,1.261419418892075,1.200320880323831,1.375915327343132,1.375915327343132,1.433163281568661,1.433163281568661,1.251083759412769)
  ,(0.4187110513297211,0.4592489584120765,16.28101926958897,0.4568601573824591,5.71690797808961,21.51699784690905,21.51699784690905,0.4187110513297211,0.4187110513297211,5.71690797808961,13.61695291249933,13.61695291249933,9.231051594470244,9.217538958776126,9.217538958776126,7.017831981914526,7.017831981914526,7.474957543618647)
  ,(1.228545138538301,0.8456359049885888,1.165396464890562,1.067715629300242,1.101232198893288,0.8047724825677131,0.8047724825677131,1.228545138538301,1.228545138538301,1.101232198893288,0.9530023407305006,0.9530023407305006,0.91721352881653,1.044849939999767,1.044849939999767,1.090773739634401,1.090773739634401,1.035549636529782)
  ,(0.8036410960048004,0.8632109905093264,3.390052849534958,0.8822485133660853,1.736919171544056,4.088990111118765,4.088990111118765,0.8036410960048004,0.8036410960048004,1.736919171544056,2.91295464133141,2.91295464133141,2.229706757724049,2.20985012622254,2.209850126222541,1.858297868668106,1.858297868668106,1.960843788679665)
  ,(1.118558773732313,0.952247942292905,1.098815420715842,1.021816348659361,0.9497741311741823,1.118018717583724,1.118018717583724,1.118558773732313,1.118558773732313,0.9497741311741823,1.033896424378953,1.033896424378953,1.006680263683604,1.06211720749674,1.06211720749674,1.076227599055633,1.076227599055633,1.043205222359721)
  ,(-0.003592677800000001,0.006681075699999994,30.27709800000001,0.005051382499999996,10.09155199999999,40.37854000000003,40.37854000000003,-0.003592677800000001,-0.003592677800000001,10.09155199999999,25.23504600000001,25.23504600000001,16.82559102523334,16.82216644073334,16.82216644073334,12.6157266611,12.6157266611,13.45922163006667)
  ,(121.14151,-0.0017802751,0.010282672,-0.0067091534,-0.012240956,0.0007090084400000001,0.0007090084400000001,121.14151,121.14151,-0.012240956,-0.00576597378,-0.00576597378,-0.004437407553333334,40.37665935081333,40.37665935081333,60.56787201311,60.56787201311,20.18862854932333)
  ,(20.194365,-0.0019288707,30.303508,-0.003354252199999999,10.089288,0.0041544129,0.0041544129,20.194365,20.194365,10.089288,5.04672120645,5.04672120645,3.3638378474,10.0959358043,10.0959358043,12.620543103225,12.620543103225,10.09767204833333)
  ,(1.2851635,0.0018918111,1.4083754,-0.0068264178,4.9112037,0.94279376,0.94279376,1.2851635,1.2851635,4.9112037,2.92699873,2.92699873,1.951963090366666,2.37972032,2.37972032,2.106081115,2.106081115,1.423766958883333)
  ,(7.263329451383068,1.718655274946576,0.9382174675479452,1.555829678690247,1.522394856945206,2.492675206161095,2.492675206161095,7.263329451383068,7.263329451383068,1.522394856945206,2.00753503155315,2.00753503155315,1.911241779350959,3.75946650482979,3.75946650482979,4.635432241468109,4.635432241468109,2.581850322612356)
  ,(4.052230811670526,0.001653716411578947,12.13835804947369,-0.005626627673684209,4.045530941157893,8.103328102157372,8.103328102157372,4.052230811670526,4.052230811670526,4.045530941157893,6.074429521657633,6.074429521657633,4.050170919908948,5.400363284995264,5.400363284995264,5.06333016666408,5.06333016666408,4.722579165532895)
  ,(7.979312787437896,0.2759162879236842,3.484395027764211,1.566428756198421,1.896432722231052,2.259136785006686,2.259136785006686,7.979312787437896,7.979312787437896,1.896432722231052,2.077784753618869,2.077784753618869,1.477161931720474,4.044960764891878,4.044960764891878,5.028548770528382,5.028548770528382,2.910270394426992)
  ,(0.4331838084551377,1.688266149648427,0.6126721767637138,1.228362198520639,0.551084634679422,1.68654151015918,1.68654151015918,0.4331838084551377,0.4331838084551377,0.551084634679422,1.118813072419301,1.118813072419301,1.30863076482901,0.8902699844312464,0.8902699844312464,0.7759984404372192,0.7759984404372192,1.033351746371087)
  ,(0.3878749146360924,1.471135685741275,0.3973658516338666,0.970124542052011,0.7411530584963391,0.05381285810825413,0.05381285810825413,0.3878749146360924,0.3878749146360924,0.7411530584963391,0.3974829583022966,0.3974829583022966,0.7553672007819562,0.3942802770802285,0.3942802770802285,0.3926789364691945,0.3926789364691945,0.6702444851113064)
  ,(1.705554643115539,1.60510479929576,0.0009603180863463678,1.44589577300558,1.042010413470588,0.002485235279404896,0.002485235279404896,1.705554643115539,1.705554643115539,1.042010413470588,0.5222478243749965,0.5222478243749965,0.8832001493485843,0.9166834306218439,0.9166834306218439,1.113901233745268,1.113901233745268,0.9670018637088696)
  ,(1.069219111067938,0.7245298929213373,1.669950404516509,0.9738373934482308,1.255596748899099,2.918949149647179,2.918949149647179,1.069219111067938,1.069219111067938,1.255596748899099,2.087272949273139,2.087272949273139,1.633025263822539,1.747921669871406,1.747921669871406,1.578246030170539,1.578246030170539,1.435347116750049)
  ,(1.794983312429255,0.5264296063171332,0.07054503366969642,1.401172564976549,0.5390105871807543,0.001806637258293289,0.001806637258293289,1.794983312429255,1.794983312429255,0.5390105871807543,0.2704086122195238,0.2704086122195238,0.3557489435853936,0.7786001789561009,0.7786001789561007,1.032695962324389,1.032695962324389,0.7223246236386135)
  ,(0.2863276994801713,1.334802064628523,1.120897010725991,1.047774481154559,1.866557827536633,0.6115437164993525,0.6115437164993525,0.2863276994801713,0.2863276994801713,1.866557827536633,1.239050772017992,1.239050772017992,1.270967869554836,0.9214764145053854,0.9214764145053854,0.762689235749082,0.762689235749082,1.044650466670871)
  ,(2.566107441077575,2.637982112588745,6.745019429020989,1.232173886747925,1.703246889088722,4.031153167192008,4.031153167192008,2.566107441077575,2.566107441077575,1.703246889088722,2.867200028140365,2.867200028140365,2.790794056289825,2.766835832452768,2.766835832452768,2.71665373460897,2.71665373460897,3.152613820952661)
  ,(3.311595942327824,0.9593792415122903,2.745522149341091,1.465823094931901,1.157603070290836,4.174194138986013,4.174194138986013,3.311595942327824,3.311595942327824,1.157603070290836,2.665898604638425,2.665898604638425,2.097058816929713,2.881131050534891,2.881131050534891,2.988747273483124,2.988747273483124,2.302352939564992)
  ,(6.911418378338329,1.273394447950587,6.040787931321557,0.3290680290204437,1.117608088513095,8.063979327355129,8.063979327355129,6.911418378338329,6.911418378338329,1.117608088513095,4.590793707934112,4.590793707934112,3.48499395460627,5.364335264735518,5.364335264735518,5.751106043136221,5.751106043136221,3.956042700416524)
)
for ($i = 0; $i -le 28; $i++) {
    $rowVals = $dataRows[$i]
    for ($j = 0; $j -le 17; $j++) {
        $ws.Cells.Item($i + 3, $j + 3).Value = $rowVals[$j]
    }
}
